# Auto-generated Excel COM-interop script
# Applies scheduled-runner market data updates to the Yojimbo_Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 429.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 500
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").ClearContents()
$ws.Range("H61").Value = 1398.6875
$ws.Range("I61").Value = 1319.2759
$ws.Range("J61").Value = 2166.3333
$ws.Range("K61").Value = 1319.2759
$ws.Range("L61").Value = 2166.3333
$ws.Range("M61").Value = -1107.2759
$ws.Range("N61").Value = -2590.3333
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("H74").Value = 1300.4634
$ws.Range("I74").Value = 797.65216
$ws.Range("J74").Value = 1942.9445
$ws.Range("K74").Value = 797.65216
$ws.Range("L74").Value = 1942.9445
$ws.Range("M74").Value = 76.34784000000002
$ws.Range("N74").Value = -3690.9445
$ws.Range("H77").Value = 1300.4634
$ws.Range("I77").Value = 797.65216
$ws.Range("J77").Value = 1942.9445
$ws.Range("K77").Value = 3988.2608
$ws.Range("L77").Value = 9714.7225
$ws.Range("M77").Value = 379.7392
$ws.Range("N77").Value = -18450.7225
$ws.Range("H97").Value = 2475.5833
$ws.Range("I97").Value = 1882.4546
$ws.Range("J97").Value = 9000
$ws.Range("K97").Value = 1882.4546
$ws.Range("L97").Value = 9000
$ws.Range("M97").Value = -1386.4546
$ws.Range("N97").Value = -9992
$ws.Range("H132").Value = 2611.8235
$ws.Range("I132").Value = 2634.6206
$ws.Range("J132").Value = 2479.6
$ws.Range("K132").Value = 7903.861800000001
$ws.Range("L132").Value = 7438.799999999999
$ws.Range("M132").Value = -5373.861800000001
$ws.Range("N132").Value = -12498.8
$ws.Range("H136").Value = 1398.6875
$ws.Range("I136").Value = 1319.2759
$ws.Range("J136").Value = 2166.3333
$ws.Range("K136").Value = 3957.8277
$ws.Range("L136").Value = 6498.999899999999
$ws.Range("M136").Value = -1407.8277
$ws.Range("N136").Value = -11598.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 14508.6
$ws.Range("I82").Value = 6815
$ws.Range("K82").Value = 6815
$ws.Range("M82").Value = -6432
$ws.Range("H85").Value = 14508.6
$ws.Range("I85").Value = 6815
$ws.Range("K85").Value = 6815
$ws.Range("M85").Value = -5489
$ws.Range("H94").Value = 827.86664
$ws.Range("I94").Value = 929.75
$ws.Range("J94").Value = 711.4286
$ws.Range("K94").Value = 929.75
$ws.Range("L94").Value = 711.4286
$ws.Range("M94").Value = -478.75
$ws.Range("N94").Value = -1613.4286
$ws.Range("H132").Value = 49390
$ws.Range("J132").Value = 49390
$ws.Range("L132").Value = 49390
$ws.Range("N132").Value = -59510
$ws.Range("H134").Value = 1122.75
$ws.Range("I134").Value = 1122.0435
$ws.Range("J134").Value = 1126
$ws.Range("K134").Value = 3366.1305
$ws.Range("L134").Value = 3378
$ws.Range("M134").Value = -831.1305000000002
$ws.Range("N134").Value = -8448

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1015.0769
$ws.Range("I22").Value = 1070.5
$ws.Range("J22").Value = 350
$ws.Range("K22").Value = 1070.5
$ws.Range("L22").Value = 350
$ws.Range("M22").Value = -720.5
$ws.Range("N22").Value = -1050

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 550.75
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 550.75
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 1652.25
$ws.Range("M92").ClearContents()
$ws.Range("N92").Value = -4148.25
$ws.Range("H105").Value = 10100
$ws.Range("J105").Value = 10100
$ws.Range("L105").Value = 30300
$ws.Range("N105").Value = -35542
$ws.Range("H137").Value = 3186.3333
$ws.Range("I137").Value = 1640
$ws.Range("J137").Value = 3669.5625
$ws.Range("K137").Value = 4920
$ws.Range("L137").Value = 11008.6875
$ws.Range("M137").Value = 180
$ws.Range("N137").Value = -21208.6875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 37.692307
$ws.Range("I2").Value = 45.5
$ws.Range("J2").Value = 34.22222
$ws.Range("K2").Value = 45.5
$ws.Range("L2").Value = 34.22222
$ws.Range("M2").Value = 67.5
$ws.Range("N2").Value = -260.22222
$ws.Range("H97").Value = 679.9231
$ws.Range("I97").Value = 502.79166
$ws.Range("K97").Value = 502.79166
$ws.Range("M97").Value = -6.791659999999979
$ws.Range("H132").Value = 2633.7
$ws.Range("I132").Value = 2450.647
$ws.Range("J132").Value = 3671
$ws.Range("K132").Value = 7351.941
$ws.Range("L132").Value = 11013
$ws.Range("M132").Value = -4821.941
$ws.Range("N132").Value = -16073
$ws.Range("H133").Value = 51682.5
$ws.Range("J133").Value = 51682.5
$ws.Range("L133").Value = 51682.5
$ws.Range("N133").Value = -61802.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 2025000
$ws.Range("I17").Value = 4000000
$ws.Range("J17").Value = 50000
$ws.Range("K17").Value = 4000000
$ws.Range("L17").Value = 50000
$ws.Range("M17").Value = -3999830
$ws.Range("N17").Value = -50340
$ws.Range("H22").Value = 348.77777
$ws.Range("I22").Value = 348.77777
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 348.77777
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -53.77776999999998
$ws.Range("N22").ClearContents()
$ws.Range("H27").Value = 348.77777
$ws.Range("I27").Value = 348.77777
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 348.77777
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -241.77777
$ws.Range("N27").ClearContents()
$ws.Range("H132").Value = 1932.875
$ws.Range("I132").Value = 1580.5227
$ws.Range("K132").Value = 4741.5681
$ws.Range("M132").Value = -2211.5681

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 14333.167
$ws.Range("I21").Value = 7499.5
$ws.Range("J21").Value = 17750
$ws.Range("K21").Value = 7499.5
$ws.Range("L21").Value = 17750
$ws.Range("M21").Value = -7264.5
$ws.Range("N21").Value = -18220
$ws.Range("H35").Value = 14333.167
$ws.Range("I35").Value = 7499.5
$ws.Range("J35").Value = 17750
$ws.Range("K35").Value = 7499.5
$ws.Range("L35").Value = 17750
$ws.Range("M35").Value = -7209.5
$ws.Range("N35").Value = -18330
$ws.Range("H100").Value = 2875.125
$ws.Range("I100").Value = 2875.125
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 5750.25
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -5209.25
$ws.Range("N100").ClearContents()
$ws.Range("H132").Value = 1231.5454
$ws.Range("I132").Value = 736.2692
$ws.Range("J132").Value = 3071.1428
$ws.Range("K132").Value = 2208.8076
$ws.Range("L132").Value = 9213.428400000001
$ws.Range("M132").Value = 321.1923999999999
$ws.Range("N132").Value = -14273.4284
